$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "price" column (D): header in D1 and a constant value (56)
# for every one of the 101 data rows (rows 2-102), mirroring the existing
# numeric "category_id" column.
$ws.Range("D1").Value = "price"
$ws.Range("D2:D102").Value = 56

# Reflect the view state left behind by the edit: scrolled near the bottom
# of the new column and with the freshly-populated D3:D102 range selected.
$excel.ActiveWindow.ScrollRow = 87
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D3:D102").Select()
